# New PO forecast model
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Weekly Quantity": append two new weekly rows (15 and 16)
# ---------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("A15").Value = 45669.99999999999
$wsWeekly.Range("B15").Value = 145
$wsWeekly.Range("A16").Value = 45676.99999999999
$wsWeekly.Range("B16").Value = 3

# ---------------------------------------------------------------
# Sheet "Monthly Trend": append one new monthly row (7)
# ---------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("A7").Value = 45688.99999999999
$wsMonthly.Range("B7").Value = 148

# ---------------------------------------------------------------
# Sheet "PO Forecast": updated forecast values + two new rows
# ---------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$wsForecast.Range("B2").Value = 91
$wsForecast.Range("B3").Value = 91
$wsForecast.Range("B4").Value = 90
$wsForecast.Range("B5").Value = 90
$wsForecast.Range("B8").Value = 83
$wsForecast.Range("B9").Value = 83
$wsForecast.Range("B11").Value = 82
$wsForecast.Range("B12").Value = 82
$wsForecast.Range("B13").Value = 81
$wsForecast.Range("B14").Value = 81

$wsForecast.Range("A15").Value = 45669.99999999999
$wsForecast.Range("B15").Value = 77
$wsForecast.Range("A16").Value = 45676.99999999999
$wsForecast.Range("B16").Value = 76
$wsForecast.Range("A17").Value = 45683.99999999999
$wsForecast.Range("B17").Value = 76
$wsForecast.Range("A18").Value = 45690.99999999999
$wsForecast.Range("B18").Value = 76
$wsForecast.Range("A19").Value = 45697.99999999999
$wsForecast.Range("B19").Value = 75
$wsForecast.Range("A20").Value = 45704.99999999999
$wsForecast.Range("B20").Value = 75
$wsForecast.Range("A21").Value = 45711.99999999999
$wsForecast.Range("B21").Value = 74
$wsForecast.Range("A22").Value = 45718.99999999999
$wsForecast.Range("B22").Value = 74

$wsForecast.Range("A23").Value = 45725.99999999999
$wsForecast.Range("B23").Value = 74
$wsForecast.Range("A24").Value = 45732.99999999999
$wsForecast.Range("B24").Value = 73

# Copy the date-column style (s="2") from an existing cell down to the
# newly added date cells so the new rows match the existing formatting.
$wsWeekly.Range("A14").Copy() | Out-Null
$wsWeekly.Range("A15:A16").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wsMonthly.Range("A6").Copy() | Out-Null
$wsMonthly.Range("A7").PasteSpecial(-4122) | Out-Null

$wsForecast.Range("A14").Copy() | Out-Null
$wsForecast.Range("A15:A24").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
